$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 4.75
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("AP5").Value = 26
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 67
$ws.Range("AW5").Value = 5.5
$ws.Range("AX5").Value = 23
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 81
$ws.Range("BA5").Value = 126
$ws.Range("G23").Value = 3.2
$ws.Range("H23").Value = 3.75
$ws.Range("L23").Value = 2.63
$ws.Range("M23").Value = 1.03
$ws.Range("N23").Value = 17
$ws.Range("O23").Value = 1.17
$ws.Range("P23").Value = 5
$ws.Range("Q23").Value = 1.57
$ws.Range("R23").Value = 2.35
$ws.Range("T23").Value = 3.5
$ws.Range("U23").Value = 1.5
$ws.Range("V23").Value = 2.5
$ws.Range("X23").Value = 19
$ws.Range("AO23").Value = 15
$ws.Range("AT23").Value = 3.5
$ws.Range("G27").Value = 2.3
$ws.Range("H27").Value = 2.92
$ws.Range("I27").Value = 3.2
$ws.Range("J27").Value = 2.95
$ws.Range("L27").Value = 3.75
$ws.Range("N27").Value = 6.1
$ws.Range("U27").Value = 1.93
$ws.Range("V27").Value = 1.7
$ws.Range("W27").Value = 6.1
$ws.Range("X27").Value = 10
$ws.Range("Z27").Value = 24
$ws.Range("AC27").Value = 6.7
$ws.Range("AG27").Value = 7.7
$ws.Range("AH27").Value = 15.5
$ws.Range("AI27").Value = 11.75
$ws.Range("AJ27").Value = 45
$ws.Range("AK27").Value = 32
$ws.Range("AN27").Value = 4
$ws.Range("AO27").Value = 12.5
$ws.Range("AQ27").Value = 55
$ws.Range("AU27").Value = 7.3
$ws.Range("AW27").Value = 4.9
$ws.Range("AX27").Value = 18
$ws.Range("AY27").Value = 27
$ws.Range("AZ27").Value = 90
